$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column K data: header year 2022 and value
$ws.Range("K4").Value = 2022
$ws.Range("K5").Value = 0.11705180708279034

# Copy style from J4/J5 onto K4/K5 so formatting matches (percent / year style)
$ws.Range("J4").Copy()
$ws.Range("K4").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("J5").Copy()
$ws.Range("K5").PasteSpecial(-4122)

# K3 should carry the thick-bottom-border style (same as J3)
$ws.Range("J3").Copy()
$ws.Range("K3").PasteSpecial(-4122)

# Columns D:K (4-11) get an explicit width of 9, replacing the old bestFit
# column J entry. (8.1 "characters" serializes to the stored width of 9.)
$ws.Range("D1:K1").ColumnWidth = 8.1

# Update selection to reflect new active cell
$ws.Range("J12").Select()
